$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear B2 and D2; update C2 and E2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.92012669354727072
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -1.0543777440885971

# Row 3: update B3:E3
$ws.Range("B3").Value = -1.3172335172622707
$ws.Range("C3").Value = -0.18180824594415326
$ws.Range("D3").Value = -2.0383289913388207
$ws.Range("E3").Value = 1.930574663407584

# Update the selection to reflect the new active range
$ws.Range("B1:E3").Select()
